$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.804.53"
$ws.Range("D3").Value = "1.953.05"
$ws.Range("E3").Value = "  +6.99%  "
$ws.Range("E4").Value = "  -0.44%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "342.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.97%  "
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4775"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.57%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4146"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +8.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.92"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08253"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.26%  "
$ws.Range("E11").Value = "  +8.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.73"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +8.16%  "
$ws.Range("D13").Value = "1.953.05"
$ws.Range("E13").Value = "  +6.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.190"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.424"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "92.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.82%  "
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("E18").Value = "  +3.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06701"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.81%  "
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("D22").Value = "29.765.07"
$ws.Range("E22").Value = "  +8.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.583"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.76%  "
$ws.Range("E24").Value = "  +4.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.263"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").Value = "2.176.65"
$ws.Range("E26").Value = "  +6.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.79%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.185"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.701"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +8.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "122.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.012"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +9.11%  "
$ws.Range("E33").Value = "  +2.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.480"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +12.75%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.689"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.25%  "
$ws.Range("E36").Value = "  +6.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06303"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.91%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02321"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.497"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.12%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.187"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.66%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6102"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "10.74"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.67%  "
$ws.Range("E43").Value = "  -0.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1894"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.33%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.389"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +32.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.262"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.26%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.54"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.58%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5707"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.81%  "
$ws.Range("E49").Value = "  +5.78%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07344"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "113.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.08%  "

Write-Host "Applied cryptos update"